$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 20: quantity, name (new part), price/piece, price all
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "Resistor 1206 SMD  620 Ohm"
$ws.Range("F20").Value = 0.02
$ws.Range("G20").Value = 0.04

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("K20").Select()
